$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (SkyBox) - count 2 -> 1, UUID -> WrapLinear, ResourceType -> Sampler, drop BindFlags (E3)
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "WrapLinear"
$ws.Range("D3").Value = "Sampler"
$ws.Range("E3").ClearContents()

# Row 4 (new) - Grid Debug / PerDebug / Buffer / ConstantBuffer / PS, slot 5
$ws.Range("A4").Value = "Grid Debug"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "PerDebug"
$ws.Range("D4").Value = "Buffer"
$ws.Range("E4").Value = '"ConstantBuffer"'
$ws.Range("G4").Value = 5

# Row 5 (was Grid Debug) - now Standard / WrapLinear / Sampler, drop BindFlags (E5), slot back to 0
$ws.Range("A5").Value = "Standard"
$ws.Range("C5").Value = "WrapLinear"
$ws.Range("D5").Value = "Sampler"
$ws.Range("E5").ClearContents()
$ws.Range("G5").Value = 0

$ws.Range("C5:H5").Select()
